# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room for the new "Bad Driver" row by inserting a blank row
#    right above the old "Totals:" row (row 8). This single insertion
#    pushes the Totals row and everything below it (the blank gap and
#    the whole "Good Drivers" table) down by exactly one row, which is
#    all that is needed to reach the new row 1:35 layout.
# ------------------------------------------------------------------
$ws.Rows(8).Insert()

# ------------------------------------------------------------------
# 2) Rewrite the "Bad Drivers" data rows (3-9) with this week's values.
# ------------------------------------------------------------------

# Row 3 (new entry, pushed to the top of the bad-driver list)
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.80.0.7"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 90.90000000000001

# Row 4 (was row 3)
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.3"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 94.7

# Row 5 (was row 4)
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6"
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = 269
$ws.Range("D5").Value = 97.59999999999999

# Row 6 - text/count unchanged, only the roaming % moved slightly
$ws.Range("D6").Value = 98.3

# Row 7 - fully unchanged, no edit needed

# Row 8 (new entry, fills the row just inserted)
$ws.Range("A8").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.0.10"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = 98.90000000000001

# Row 9 - Totals (was row 8 before the insert); update the rollup values
$ws.Range("B9").Value = 20
$ws.Range("C9").Value = 339

# ------------------------------------------------------------------
# 3) A handful of sample counts in the "Good Drivers" table changed
#    for this week's refresh (the rows themselves already shifted
#    down by one thanks to the insert above, carrying their text /
#    date-vintage cells with them intact).
# ------------------------------------------------------------------

# "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4" -> now row 19
$ws.Range("B19").Value = 449371

# "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1" -> now row 20
$ws.Range("B20").Value = 14968

# "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9" -> now row 25
$ws.Range("B25").Value = 77999
